$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2022 column (M) alongside the existing year columns (D:L = 2013..2021).
$ws.Range("M4").Value = 2022

$yearValues = @{
    5  = 24.6
    6  = 40.700000000000003
    7  = 20.7
    8  = 26.6
    9  = 44.5
    10 = 21.9
    11 = 21.9
    12 = 35.299999999999997
    13 = 17.600000000000001
    14 = 28
    15 = 44.9
    16 = 21.5
    17 = 36.200000000000003
    18 = 53.1
    19 = 33.4
    20 = 20.2
    21 = 15.4
    22 = 20.5
    23 = 27.1
    24 = 36.1
    25 = 25.2
    26 = 24.2
    27 = 46.5
    28 = 20.3
    29 = 40.5
    30 = 44.5
}

# M4's formatting matches K4/L4 (the year header row).
$ws.Range("K4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 2022

foreach ($row in 5..29) {
    $ws.Range("M$row").Value = $yearValues[$row]
    $ws.Range("K$row").Copy()
    $ws.Range("M$row").PasteSpecial(-4122)
    $ws.Range("M$row").Value = $yearValues[$row]
}

# Row 14 gets an explicit "0.0" number format (new style), unlike the other rows.
$ws.Range("M14").NumberFormat = "0.0"

# Row 30 (bottom, thick border) copies its style from L30 instead of K30.
$ws.Range("M30").Value = $yearValues[30]
$ws.Range("L30").Copy()
$ws.Range("M30").PasteSpecial(-4122)
$ws.Range("M30").Value = $yearValues[30]

$excel.CutCopyMode = $false

$ws.Range("N7").Select()
